# Apply the authored content edits to 程序设计ppt.pptx
#
# 1) Title slide (slide 1): the credits textbox ("文本框 36", shape id 37)
#    originally reads "    成员：游洋 呙自桥 成俊宏" across several runs.
#    The author trimmed the roster down to just the "成员：" label.
# 2) Slide 5 ("通讯录"/THREE content slide): the 5th bullet's text
#    ("排序与" + "优化") was re-typed so the two runs collapse into a
#    single run reading "排序与优化" (no visible text change, just the
#    run split going away).

$p = $ppt.ActivePresentation

# --- Slide 1: trim the member list down to "成员：" -----------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Id -eq 37) {
        $shp.TextFrame.TextRange.Text = "    成员："
    }
}

# --- Slide 5: re-set the "排序与优化" bullet so the runs merge -------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.Id -eq 4) {
        $tr = $shp.TextFrame.TextRange
        $fullText = $tr.Text
        $idx = $fullText.IndexOf("排序与优化")
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, 5)
            $sub.Text = "排序与优化"
        }
    }
}
